$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 539075.75
$ws.Range("I15").Value = 539075.75
$ws.Range("K15").Value = 1617227.25
$ws.Range("M15").Value = -1617058.25
$ws.Range("H18").Value = 200200270
$ws.Range("I18").Value = 200200270
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 200200270
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -200199986
$ws.Range("N18").ClearContents()
$ws.Range("H98").Value = 1459.4736
$ws.Range("I98").Value = 1460.5883
$ws.Range("K98").Value = 1460.5883
$ws.Range("M98").Value = 37.41170000000011
$ws.Range("H112").Value = 1536.2
$ws.Range("I112").Value = 250
$ws.Range("J112").Value = 1648.0435
$ws.Range("K112").Value = 750
$ws.Range("L112").Value = 4944.1305
$ws.Range("M112").Value = 358
$ws.Range("N112").Value = -7160.1305
$ws.Range("H122").Value = 1459.4736
$ws.Range("I122").Value = 1460.5883
$ws.Range("K122").Value = 4381.7649
$ws.Range("M122").Value = -1931.7649
$ws.Range("H124").Value = 35543.6
$ws.Range("J124").Value = 35543.6
$ws.Range("L124").Value = 35543.6
$ws.Range("N124").Value = -45363.6

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 476.25
$ws.Range("I5").Value = 402.5
$ws.Range("J5").Value = 550
$ws.Range("K5").Value = 402.5
$ws.Range("L5").Value = 550
$ws.Range("M5").Value = -290.5
$ws.Range("N5").Value = -774
$ws.Range("H32").Value = 18871980
$ws.Range("I32").Value = 21741952
$ws.Range("J32").Value = 12154.857
$ws.Range("K32").Value = 21741952
$ws.Range("L32").Value = 12154.857
$ws.Range("M32").Value = -21741665
$ws.Range("N32").Value = -12728.857
$ws.Range("H122").Value = 1451.8334
$ws.Range("I122").Value = 1474.7273
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 4424.1819
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -1974.1819
$ws.Range("N122").Value = -8500

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 476.25
$ws.Range("I4").Value = 402.5
$ws.Range("J4").Value = 550
$ws.Range("K4").Value = 402.5
$ws.Range("L4").Value = 550
$ws.Range("M4").Value = -287.5
$ws.Range("N4").Value = -780
$ws.Range("H12").Value = 229.85715
$ws.Range("I12").Value = 229.85715
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 229.85715
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -61.85714999999999
$ws.Range("N12").ClearContents()
$ws.Range("H86").Value = 2742.889
$ws.Range("I86").Value = 2216.8333
$ws.Range("J86").Value = 3795
$ws.Range("K86").Value = 2216.8333
$ws.Range("L86").Value = 3795
$ws.Range("M86").Value = -1093.8333
$ws.Range("N86").Value = -6041
$ws.Range("H89").Value = 2742.889
$ws.Range("I89").Value = 2216.8333
$ws.Range("J89").Value = 3795
$ws.Range("K89").Value = 11084.1665
$ws.Range("L89").Value = 18975
$ws.Range("M89").Value = -5468.166499999999
$ws.Range("N89").Value = -30207
$ws.Range("H105").Value = 2476.71
$ws.Range("I105").Value = 1388.7084
$ws.Range("J105").Value = 2820.2896
$ws.Range("K105").Value = 1388.7084
$ws.Range("L105").Value = 2820.2896
$ws.Range("M105").Value = 358.2916
$ws.Range("N105").Value = -6314.2896

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 84.117645
$ws.Range("I7").Value = 77.333336
$ws.Range("K7").Value = 77.333336
$ws.Range("M7").Value = 35.666664
$ws.Range("H11").Value = 9905
$ws.Range("I11").Value = 310
$ws.Range("J11").Value = 19500
$ws.Range("K11").Value = 310
$ws.Range("L11").Value = 19500
$ws.Range("M11").Value = -170
$ws.Range("N11").Value = -19780
$ws.Range("H31").Value = 1719.7872
$ws.Range("I31").Value = 1297.1538
$ws.Range("J31").Value = 3780.125
$ws.Range("K31").Value = 1297.1538
$ws.Range("L31").Value = 3780.125
$ws.Range("M31").Value = -1002.1538
$ws.Range("N31").Value = -4370.125
$ws.Range("H34").Value = 1719.7872
$ws.Range("I34").Value = 1297.1538
$ws.Range("J34").Value = 3780.125
$ws.Range("K34").Value = 1297.1538
$ws.Range("L34").Value = 3780.125
$ws.Range("M34").Value = -1095.1538
$ws.Range("N34").Value = -4184.125
$ws.Range("H86").Value = 4272.357
$ws.Range("I86").Value = 5872.143
$ws.Range("J86").Value = 2672.5715
$ws.Range("K86").Value = 5872.143
$ws.Range("L86").Value = 2672.5715
$ws.Range("M86").Value = -4749.143
$ws.Range("N86").Value = -4918.5715
$ws.Range("H89").Value = 4272.357
$ws.Range("I89").Value = 5872.143
$ws.Range("J89").Value = 2672.5715
$ws.Range("K89").Value = 29360.715
$ws.Range("L89").Value = 13362.8575
$ws.Range("M89").Value = -23744.715
$ws.Range("N89").Value = -24594.8575

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1452.2941
$ws.Range("I17").Value = 419.8
$ws.Range("J17").Value = 1882.5
$ws.Range("K17").Value = 1259.4
$ws.Range("L17").Value = 5647.5
$ws.Range("M17").Value = -1090.4
$ws.Range("N17").Value = -5985.5
$ws.Range("H39").Value = 3657.3333
$ws.Range("I39").Value = 1500
$ws.Range("J39").Value = 3927
$ws.Range("K39").Value = 4500
$ws.Range("L39").Value = 11781
$ws.Range("M39").Value = -4206
$ws.Range("N39").Value = -12369
$ws.Range("H68").Value = 536.25
$ws.Range("I68").Value = 372.5
$ws.Range("J68").Value = 700
$ws.Range("K68").Value = 1117.5
$ws.Range("L68").Value = 2100
$ws.Range("M68").Value = -306.5
$ws.Range("N68").Value = -3722
$ws.Range("H71").Value = 536.25
$ws.Range("I71").Value = 372.5
$ws.Range("J71").Value = 700
$ws.Range("K71").Value = 3352.5
$ws.Range("L71").Value = 6300
$ws.Range("M71").Value = 703.5
$ws.Range("N71").Value = -14412
$ws.Range("H98").Value = 429.42856
$ws.Range("I98").Value = 334.33334
$ws.Range("K98").Value = 1003.00002
$ws.Range("M98").Value = 494.9999799999999
$ws.Range("H118").Value = 2092.9375
$ws.Range("I118").Value = 641
$ws.Range("K118").Value = 1923
$ws.Range("M118").Value = -680
$ws.Range("H122").Value = 5641.905
$ws.Range("I122").Value = 9734.739
$ws.Range("J122").Value = 687.4211
$ws.Range("K122").Value = 87612.651
$ws.Range("L122").Value = 6186.7899
$ws.Range("M122").Value = -85162.651
$ws.Range("N122").Value = -11086.7899
$ws.Range("H125").Value = 4155.385
$ws.Range("J125").Value = 4418.3335
$ws.Range("L125").Value = 13255.0005
$ws.Range("N125").Value = -23095.0005

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5872.3335
$ws.Range("I70").Value = 5934.643
$ws.Range("K70").Value = 5934.643
$ws.Range("M70").Value = -5664.643
$ws.Range("H73").Value = 5872.3335
$ws.Range("I73").Value = 5934.643
$ws.Range("K73").Value = 5934.643
$ws.Range("M73").Value = -4998.643
$ws.Range("H102").Value = 1657.4517
$ws.Range("I102").Value = 1705
$ws.Range("J102").Value = 1520.75
$ws.Range("K102").Value = 1705
$ws.Range("L102").Value = 1520.75
$ws.Range("M102").Value = -83
$ws.Range("N102").Value = -4764.75

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 2000
$ws.Range("I20").Value = 1000
$ws.Range("J20").Value = 3000
$ws.Range("K20").Value = 1000
$ws.Range("L20").Value = 3000
$ws.Range("M20").Value = -774
$ws.Range("N20").Value = -3452

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 630.0714
$ws.Range("I107").Value = 639.8461
$ws.Range("J107").Value = 503
$ws.Range("K107").Value = 1919.5383
$ws.Range("L107").Value = 1509
$ws.Range("M107").Value = 0.4617000000000644
$ws.Range("N107").Value = -5349
$ws.Range("H132").Value = 2456.0334
$ws.Range("I132").Value = 2157.2083
$ws.Range("K132").Value = 6471.624899999999
$ws.Range("M132").Value = -3941.624899999999
$ws.Range("H136").Value = 3424.94
$ws.Range("I136").Value = 3431.55
$ws.Range("J136").Value = 3398.5
$ws.Range("K136").Value = 10294.65
$ws.Range("L136").Value = 10195.5
$ws.Range("M136").Value = -7744.650000000001
$ws.Range("N136").Value = -15295.5
